$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.130952380952381
$ws.Range("C2").Value = 0.469047619047619
$ws.Range("D2").Value = 0.5435714285714285
$ws.Range("E2").Value = 0.7372729674763808
$ws.Range("F2").Value = 0.7343449363085427
$ws.Range("G2").Value = 42

$ws.Range("B3").Value = 0.006862445399837091
$ws.Range("C3").Value = 0.6039302007301544
$ws.Range("D3").Value = 0.7234397629857844
$ws.Range("E3").Value = 0.8505526221144606
$ws.Range("F3").Value = 0.8610908701029418
$ws.Range("G3").Value = 41

$ws.Range("B4").Value = 0.1524205693483774
$ws.Range("C4").Value = 0.6423595844912023
$ws.Range("D4").Value = 0.8775854164232644
$ws.Range("E4").Value = 0.9367952905642003
$ws.Range("F4").Value = 0.9360875405385815
$ws.Range("G4").Value = 40

$ws.Range("B5").Value = 0.04980756042726093
$ws.Range("C5").Value = 0.7196857601028602
$ws.Range("D5").Value = 0.9118160530585012
$ws.Range("E5").Value = 0.9548905974290988
$ws.Range("F5").Value = 0.9660564865707858
$ws.Range("G5").Value = 39

$ws.Range("B6").Value = 0.1685415784894023
$ws.Range("C6").Value = 0.7315745272820626
$ws.Range("D6").Value = 0.9669495865848617
$ws.Range("E6").Value = 0.983335947977527
$ws.Range("F6").Value = 0.9817888564550711
$ws.Range("G6").Value = 38

$ws.Range("B7").Value = 0.05236864345981531
$ws.Range("C7").Value = 0.6876603977094401
$ws.Range("D7").Value = 0.8135759198832745
$ws.Range("E7").Value = 0.9019844343907907
$ws.Range("F7").Value = 0.9128836707474132
$ws.Range("G7").Value = 37

$ws.Range("B8").Value = 0.1715962988235919
$ws.Range("C8").Value = 0.6741482612662726
$ws.Range("D8").Value = 0.9067807233488873
$ws.Range("E8").Value = 0.9522503469933142
$ws.Range("F8").Value = 0.9499485040004092
$ws.Range("G8").Value = 36

$ws.Range("B9").Value = 0.1127370711831118
$ws.Range("C9").Value = 0.5725331809908775
$ws.Range("D9").Value = 0.6431051428235416
$ws.Range("E9").Value = 0.801938365975554
$ws.Range("F9").Value = 0.805565974699134
$ws.Range("G9").Value = 35

$ws.Range("B10").Value = 0.1706140055022492
$ws.Range("C10").Value = 0.666628778494778
$ws.Range("D10").Value = 0.750579926290491
$ws.Range("E10").Value = 0.8663601596856189
$ws.Range("F10").Value = 0.8621679294376572
$ws.Range("G10").Value = 34

$ws.Range("B11").Value = 0.2039437044841572
$ws.Range("C11").Value = 0.6558282946641647
$ws.Range("D11").Value = 0.7409461765722168
$ws.Range("E11").Value = 0.8607823049832152
$ws.Range("F11").Value = 0.8492396173402222
$ws.Range("G11").Value = 33

